$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "60.706.05"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -1.42%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.384.49"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -1.88%  "

$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.02%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "568.08"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -2.22%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "141.27"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -3.26%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.384.95"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -1.90%  "

$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -0.26%  "

$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -1.97%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.123"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -1.66%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.398"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +1.91%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "3.962.04"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -1.88%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "28.34"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +1.02%  "

$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +1.58%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0000170"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -2.01%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.380.82"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -2.02%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "60.819.00"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -1.40%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.23"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -0.04%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "14.04"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -2.04%  "

$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -5.43%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "383.49"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -1.44%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.561"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -0.94%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "73.63"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -0.05%  "

$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +0.67%  "

$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -5.12%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "3.522.38"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -2.00%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.179"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -1.66%  "

$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -0.04%  "

$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -2.94%  "

$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -2.32%  "

$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -2.19%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.42"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -3.44%  "

$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -0.01%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "23.71"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -1.65%  "

$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -0.42%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "166.16"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -0.50%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "5.02"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -2.37%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.416.03"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -1.76%  "

$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -4.72%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "27.78"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +1.08%  "

$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -1.14%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.781"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -3.01%  "

$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -0.07%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "41.82"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -1.60%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "4.41"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -2.32%  "

$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -2.92%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.514.76"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -2.18%  "

$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -3.53%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "23.66"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +2.87%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "6.84"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -1.70%  "

